# branch_wise_stock_status.xlsx — refresh the per-branch stock status
# figures (columns D:H i.e. the stock-count buckets) for rows 2-32 with
# the latest counts. Column A (branch), B (total) and C stay unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 33  # D2: 34 -> 33
$ws.Cells.Item(2, 6).Value = 9  # F2: 13 -> 9
$ws.Cells.Item(2, 7).Value = 14  # G2: 12 -> 14
$ws.Cells.Item(2, 8).Value = 40  # H2: 37 -> 40
# Row 3
$ws.Cells.Item(3, 4).Value = 37  # D3: 35 -> 37
$ws.Cells.Item(3, 5).Value = 16  # E3: 14 -> 16
$ws.Cells.Item(3, 6).Value = 9  # F3: 13 -> 9
$ws.Cells.Item(3, 7).Value = 8  # G3: 9 -> 8
$ws.Cells.Item(3, 8).Value = 36  # H3: 35 -> 36
# Row 4
$ws.Cells.Item(4, 4).Value = 35  # D4: 34 -> 35
$ws.Cells.Item(4, 6).Value = 6  # F4: 8 -> 6
$ws.Cells.Item(4, 8).Value = 32  # H4: 31 -> 32
# Row 5
$ws.Cells.Item(5, 4).Value = 34  # D5: 35 -> 34
$ws.Cells.Item(5, 5).Value = 14  # E5: 13 -> 14
# Row 6
$ws.Cells.Item(6, 4).Value = 32  # D6: 34 -> 32
$ws.Cells.Item(6, 5).Value = 12  # E6: 11 -> 12
$ws.Cells.Item(6, 6).Value = 3  # F6: 4 -> 3
$ws.Cells.Item(6, 7).Value = 18  # G6: 19 -> 18
$ws.Cells.Item(6, 8).Value = 41  # H6: 38 -> 41
# Row 7
$ws.Cells.Item(7, 5).Value = 11  # E7: 17 -> 11
$ws.Cells.Item(7, 6).Value = 9  # F7: 7 -> 9
$ws.Cells.Item(7, 7).Value = 10  # G7: 11 -> 10
$ws.Cells.Item(7, 8).Value = 46  # H7: 41 -> 46
# Row 8
$ws.Cells.Item(8, 4).Value = 32  # D8: 34 -> 32
$ws.Cells.Item(8, 6).Value = 6  # F8: 4 -> 6
$ws.Cells.Item(8, 7).Value = 18  # G8: 13 -> 18
$ws.Cells.Item(8, 8).Value = 35  # H8: 40 -> 35
# Row 9
$ws.Cells.Item(9, 4).Value = 31  # D9: 28 -> 31
$ws.Cells.Item(9, 5).Value = 7  # E9: 8 -> 7
$ws.Cells.Item(9, 8).Value = 52  # H9: 54 -> 52
# Row 10
$ws.Cells.Item(10, 4).Value = 39  # D10: 40 -> 39
$ws.Cells.Item(10, 6).Value = 8  # F10: 6 -> 8
$ws.Cells.Item(10, 7).Value = 13  # G10: 12 -> 13
$ws.Cells.Item(10, 8).Value = 36  # H10: 38 -> 36
# Row 11
$ws.Cells.Item(11, 4).Value = 31  # D11: 32 -> 31
$ws.Cells.Item(11, 6).Value = 4  # F11: 10 -> 4
$ws.Cells.Item(11, 7).Value = 8  # G11: 7 -> 8
$ws.Cells.Item(11, 8).Value = 49  # H11: 43 -> 49
# Row 12
$ws.Cells.Item(12, 4).Value = 33  # D12: 34 -> 33
$ws.Cells.Item(12, 5).Value = 15  # E12: 16 -> 15
$ws.Cells.Item(12, 6).Value = 9  # F12: 10 -> 9
$ws.Cells.Item(12, 7).Value = 6  # G12: 5 -> 6
$ws.Cells.Item(12, 8).Value = 43  # H12: 41 -> 43
# Row 13
$ws.Cells.Item(13, 4).Value = 33  # D13: 34 -> 33
$ws.Cells.Item(13, 5).Value = 25  # E13: 23 -> 25
$ws.Cells.Item(13, 6).Value = 8  # F13: 9 -> 8
$ws.Cells.Item(13, 7).Value = 13  # G13: 15 -> 13
$ws.Cells.Item(13, 8).Value = 27  # H13: 25 -> 27
# Row 14
$ws.Cells.Item(14, 4).Value = 35  # D14: 34 -> 35
$ws.Cells.Item(14, 5).Value = 12  # E14: 11 -> 12
$ws.Cells.Item(14, 6).Value = 6  # F14: 9 -> 6
$ws.Cells.Item(14, 7).Value = 10  # G14: 9 -> 10
# Row 15
$ws.Cells.Item(15, 4).Value = 29  # D15: 30 -> 29
$ws.Cells.Item(15, 5).Value = 8  # E15: 10 -> 8
$ws.Cells.Item(15, 6).Value = 7  # F15: 6 -> 7
$ws.Cells.Item(15, 8).Value = 52  # H15: 50 -> 52
# Row 16
$ws.Cells.Item(16, 4).Value = 34  # D16: 36 -> 34
$ws.Cells.Item(16, 5).Value = 16  # E16: 13 -> 16
$ws.Cells.Item(16, 6).Value = 5  # F16: 7 -> 5
$ws.Cells.Item(16, 7).Value = 17  # G16: 14 -> 17
$ws.Cells.Item(16, 8).Value = 34  # H16: 36 -> 34
# Row 17
$ws.Cells.Item(17, 4).Value = 38  # D17: 36 -> 38
$ws.Cells.Item(17, 5).Value = 20  # E17: 22 -> 20
$ws.Cells.Item(17, 6).Value = 9  # F17: 10 -> 9
$ws.Cells.Item(17, 7).Value = 7  # G17: 8 -> 7
$ws.Cells.Item(17, 8).Value = 32  # H17: 30 -> 32
# Row 18
$ws.Cells.Item(18, 5).Value = 7  # E18: 8 -> 7
$ws.Cells.Item(18, 7).Value = 1  # G18: 4 -> 1
$ws.Cells.Item(18, 8).Value = 65  # H18: 61 -> 65
# Row 19
$ws.Cells.Item(19, 4).Value = 32  # D19: 30 -> 32
$ws.Cells.Item(19, 5).Value = 9  # E19: 11 -> 9
$ws.Cells.Item(19, 7).Value = 8  # G19: 6 -> 8
$ws.Cells.Item(19, 8).Value = 54  # H19: 56 -> 54
# Row 20
$ws.Cells.Item(20, 4).Value = 35  # D20: 36 -> 35
$ws.Cells.Item(20, 5).Value = 11  # E20: 13 -> 11
$ws.Cells.Item(20, 6).Value = 11  # F20: 10 -> 11
$ws.Cells.Item(20, 7).Value = 16  # G20: 13 -> 16
$ws.Cells.Item(20, 8).Value = 33  # H20: 34 -> 33
# Row 21
$ws.Cells.Item(21, 5).Value = 15  # E21: 16 -> 15
$ws.Cells.Item(21, 7).Value = 7  # G21: 5 -> 7
$ws.Cells.Item(21, 8).Value = 46  # H21: 47 -> 46
# Row 22
$ws.Cells.Item(22, 4).Value = 35  # D22: 36 -> 35
$ws.Cells.Item(22, 5).Value = 19  # E22: 20 -> 19
$ws.Cells.Item(22, 6).Value = 6  # F22: 5 -> 6
$ws.Cells.Item(22, 7).Value = 10  # G22: 11 -> 10
$ws.Cells.Item(22, 8).Value = 36  # H22: 34 -> 36
# Row 23
$ws.Cells.Item(23, 5).Value = 15  # E23: 12 -> 15
$ws.Cells.Item(23, 6).Value = 12  # F23: 13 -> 12
$ws.Cells.Item(23, 7).Value = 10  # G23: 12 -> 10
# Row 24
$ws.Cells.Item(24, 4).Value = 39  # D24: 37 -> 39
$ws.Cells.Item(24, 5).Value = 16  # E24: 20 -> 16
$ws.Cells.Item(24, 6).Value = 8  # F24: 6 -> 8
# Row 25
$ws.Cells.Item(25, 6).Value = 4  # F25: 3 -> 4
$ws.Cells.Item(25, 7).Value = 11  # G25: 10 -> 11
$ws.Cells.Item(25, 8).Value = 56  # H25: 58 -> 56
# Row 26
$ws.Cells.Item(26, 7).Value = 9  # G26: 12 -> 9
$ws.Cells.Item(26, 8).Value = 44  # H26: 41 -> 44
# Row 27
$ws.Cells.Item(27, 4).Value = 29  # D27: 30 -> 29
$ws.Cells.Item(27, 5).Value = 13  # E27: 14 -> 13
$ws.Cells.Item(27, 8).Value = 48  # H27: 46 -> 48
# Row 28
$ws.Cells.Item(28, 4).Value = 40  # D28: 38 -> 40
$ws.Cells.Item(28, 5).Value = 12  # E28: 13 -> 12
$ws.Cells.Item(28, 6).Value = 7  # F28: 5 -> 7
$ws.Cells.Item(28, 7).Value = 8  # G28: 10 -> 8
$ws.Cells.Item(28, 8).Value = 39  # H28: 40 -> 39
# Row 29
$ws.Cells.Item(29, 4).Value = 33  # D29: 32 -> 33
$ws.Cells.Item(29, 5).Value = 7  # E29: 11 -> 7
$ws.Cells.Item(29, 6).Value = 6  # F29: 7 -> 6
$ws.Cells.Item(29, 7).Value = 11  # G29: 9 -> 11
$ws.Cells.Item(29, 8).Value = 49  # H29: 47 -> 49
# Row 30
$ws.Cells.Item(30, 4).Value = 36  # D30: 34 -> 36
$ws.Cells.Item(30, 5).Value = 12  # E30: 15 -> 12
$ws.Cells.Item(30, 6).Value = 14  # F30: 11 -> 14
$ws.Cells.Item(30, 7).Value = 13  # G30: 16 -> 13
$ws.Cells.Item(30, 8).Value = 31  # H30: 30 -> 31
# Row 31
$ws.Cells.Item(31, 4).Value = 32  # D31: 33 -> 32
$ws.Cells.Item(31, 5).Value = 24  # E31: 23 -> 24
$ws.Cells.Item(31, 6).Value = 5  # F31: 6 -> 5
$ws.Cells.Item(31, 7).Value = 12  # G31: 13 -> 12
$ws.Cells.Item(31, 8).Value = 33  # H31: 31 -> 33
# Row 32
$ws.Cells.Item(32, 4).Value = 33  # D32: 36 -> 33
$ws.Cells.Item(32, 7).Value = 10  # G32: 7 -> 10
